$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) onto the new
# header cells so they pick up the same bold/centered/bordered format.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 100
    $ws.Cells.Item($row, 31).Value = 62
    $ws.Cells.Item($row, 32).Value = 0
}
